$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 231.18182
$ws.Range("I9").Value = 194.8
$ws.Range("J9").Value = 261.5
$ws.Range("K9").Value = 194.8
$ws.Range("L9").Value = 261.5
$ws.Range("M9").Value = -25.80000000000001
$ws.Range("N9").Value = -599.5
$ws.Range("H12").Value = 340
$ws.Range("I12").Value = 340
$ws.Range("K12").Value = 340
$ws.Range("M12").Value = -170
$ws.Range("H21").Value = 19166.334
$ws.Range("J21").Value = 18999
$ws.Range("L21").Value = 18999
$ws.Range("N21").Value = -19935
$ws.Range("H23").Value = 19166.334
$ws.Range("J23").Value = 18999
$ws.Range("L23").Value = 18999
$ws.Range("N23").Value = -19467
$ws.Range("H32").Value = 2833.3333
$ws.Range("J32").Value = 3000
$ws.Range("L32").Value = 3000
$ws.Range("N32").Value = -3652
$ws.Range("H40").Value = 1587.25
$ws.Range("I40").Value = 1528.4286
$ws.Range("J40").Value = 1999
$ws.Range("K40").Value = 1528.4286
$ws.Range("L40").Value = 1999
$ws.Range("M40").Value = -1353.4286
$ws.Range("N40").Value = -2349
$ws.Range("H41").Value = 532.6667
$ws.Range("I41").Value = 532.6667
$ws.Range("K41").Value = 532.6667
$ws.Range("M41").Value = -92.66669999999999

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1158.25
$ws.Range("I2").Value = 877.6667
$ws.Range("K2").Value = 877.6667
$ws.Range("M2").Value = -764.6667
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H94").Value = 47249.5
$ws.Range("J94").Value = 47249.5
$ws.Range("L94").Value = 47249.5
$ws.Range("N94").Value = -49051.5
$ws.Range("H116").Value = 1158.25
$ws.Range("I116").Value = 877.6667
$ws.Range("K116").Value = 877.6667
$ws.Range("M116").Value = 1416.3333
$ws.Range("H122").Value = 12375
$ws.Range("J122").Value = 9500
$ws.Range("L122").Value = 28500
$ws.Range("N122").Value = -33400

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1158.25
$ws.Range("I3").Value = 877.6667
$ws.Range("K3").Value = 877.6667
$ws.Range("M3").Value = -763.6667
$ws.Range("H107").Value = 1750.4546
$ws.Range("I107").Value = 1750.4546
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 1750.4546
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 169.5454
$ws.Range("N107").ClearContents()
$ws.Range("H134").Value = 3785.2354
$ws.Range("I134").Value = 3614.6
$ws.Range("K134").Value = 10843.8
$ws.Range("M134").Value = -8308.799999999999

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 60.633335
$ws.Range("I7").Value = 74.26316
$ws.Range("J7").Value = 37.090908
$ws.Range("K7").Value = 74.26316
$ws.Range("L7").Value = 37.090908
$ws.Range("M7").Value = 38.73684
$ws.Range("N7").Value = -263.090908
$ws.Range("H10").Value = 1405.4546
$ws.Range("I10").Value = 1055.5
$ws.Range("K10").Value = 1055.5
$ws.Range("M10").Value = -916.5
$ws.Range("H50").Value = 42947.5
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 42947.5
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 42947.5
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -44197.5
$ws.Range("H60").Value = 17770.715
$ws.Range("I60").Value = 13583.333
$ws.Range("J60").Value = 42895
$ws.Range("K60").Value = 13583.333
$ws.Range("L60").Value = 42895
$ws.Range("M60").Value = -13072.333
$ws.Range("N60").Value = -43917
$ws.Range("H82").Value = 65000
$ws.Range("J82").Value = 65000
$ws.Range("L82").Value = 65000
$ws.Range("N82").Value = -65722
$ws.Range("H85").Value = 65000
$ws.Range("J85").Value = 65000
$ws.Range("L85").Value = 65000
$ws.Range("N85").Value = -67496
$ws.Range("H106").Value = 95492.75
$ws.Range("J106").Value = 95492.75
$ws.Range("L106").Value = 95492.75
$ws.Range("N106").Value = -98016.75
$ws.Range("H111").Value = 99000
$ws.Range("J111").Value = 99000
$ws.Range("L111").Value = 99000
$ws.Range("N111").Value = -107180
$ws.Range("H122").Value = 10387
$ws.Range("I122").Value = 728.6667
$ws.Range("J122").Value = 24874.5
$ws.Range("K122").Value = 2186.0001
$ws.Range("L122").Value = 74623.5
$ws.Range("M122").Value = 263.9998999999998
$ws.Range("N122").Value = -79523.5
$ws.Range("H132").Value = 6697.6665
$ws.Range("I132").Value = 6132.3335
$ws.Range("K132").Value = 18397.0005
$ws.Range("M132").Value = -15867.0005
$ws.Range("H134").Value = 1105.25
$ws.Range("I134").Value = 1070.3334
$ws.Range("K134").Value = 3211.0002
$ws.Range("M134").Value = -676.0001999999999

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 2989.8462
$ws.Range("I14").Value = 2989.8462
$ws.Range("K14").Value = 8969.5386
$ws.Range("M14").Value = -8796.5386
$ws.Range("H23").Value = 67.73333
$ws.Range("I23").Value = 32.77778
$ws.Range("J23").Value = 120.166664
$ws.Range("K23").Value = 98.33333999999999
$ws.Range("L23").Value = 360.499992
$ws.Range("M23").Value = 136.66666
$ws.Range("N23").Value = -830.499992
$ws.Range("H86").Value = 351.2
$ws.Range("I86").Value = 452.33334
$ws.Range("J86").Value = 199.5
$ws.Range("K86").Value = 1357.00002
$ws.Range("L86").Value = 598.5
$ws.Range("M86").Value = -171.0000199999999
$ws.Range("N86").Value = -2970.5
$ws.Range("H89").Value = 351.2
$ws.Range("I89").Value = 452.33334
$ws.Range("J89").Value = 199.5
$ws.Range("K89").Value = 4071.00006
$ws.Range("L89").Value = 1795.5
$ws.Range("M89").Value = 1856.99994
$ws.Range("N89").Value = -13651.5

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4149.6665
$ws.Range("I80").Value = 4974.5
$ws.Range("J80").Value = 2500
$ws.Range("K80").Value = 4974.5
$ws.Range("L80").Value = 2500
$ws.Range("M80").Value = -3976.5
$ws.Range("N80").Value = -4496
$ws.Range("H83").Value = 4149.6665
$ws.Range("I83").Value = 4974.5
$ws.Range("J83").Value = 2500
$ws.Range("K83").Value = 24872.5
$ws.Range("L83").Value = 12500
$ws.Range("M83").Value = -19880.5
$ws.Range("N83").Value = -22484
$ws.Range("H99").Value = 18999.334
$ws.Range("I99").Value = 3499.5
$ws.Range("J99").Value = 49999
$ws.Range("K99").Value = 3499.5
$ws.Range("L99").Value = 49999
$ws.Range("M99").Value = -1253.5
$ws.Range("N99").Value = -54491
$ws.Range("H126").Value = 7685
$ws.Range("J126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("N126").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1850
$ws.Range("I22").Value = 1700
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1700
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1405
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1850
$ws.Range("I27").Value = 1700
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1700
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1593
$ws.Range("N27").Value = -2214
$ws.Range("H46").Value = 492
$ws.Range("I46").Value = 492
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 492
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -304
$ws.Range("N46").ClearContents()
$ws.Range("H61").Value = 4390.6
$ws.Range("I61").Value = 3738.25
$ws.Range("K61").Value = 3738.25
$ws.Range("M61").Value = -3536.25
$ws.Range("H93").Value = 1730.875
$ws.Range("I93").Value = 1612
$ws.Range("K93").Value = 1612
$ws.Range("M93").Value = -364
$ws.Range("H113").Value = 4390.6
$ws.Range("I113").Value = 3738.25
$ws.Range("K113").Value = 3738.25
$ws.Range("M113").Value = -1568.25
$ws.Range("H122").Value = 3438.3845
$ws.Range("I122").Value = 3180
$ws.Range("J122").Value = 3599.875
$ws.Range("K122").Value = 9540
$ws.Range("L122").Value = 10799.625
$ws.Range("M122").Value = -7090
$ws.Range("N122").Value = -15699.625

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 1114221.9
$ws.Range("J14").Value = 3499.625
$ws.Range("L14").Value = 3499.625
$ws.Range("N14").Value = -3835.625
$ws.Range("H16").Value = 42000
$ws.Range("J16").Value = 42000
$ws.Range("L16").Value = 42000
$ws.Range("N16").Value = -42584
$ws.Range("H107").Value = 1984
$ws.Range("I107").Value = 0
$ws.Range("J107").Value = 1984
$ws.Range("K107").Value = 0
$ws.Range("L107").Value = 5952
$ws.Range("M107").ClearContents()
$ws.Range("N107").Value = -9792
$ws.Range("H132").Value = 1375.4286
$ws.Range("I132").Value = 910.75
$ws.Range("K132").Value = 2732.25
$ws.Range("M132").Value = -202.25
